$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to be treated as plain text so values like
# "352.17" or "0.0453" aren't reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Rows 43 and 44 swap their Coin/Link content (and get new Price/Volume values)
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").Value = "120.30"
$ws.Range("E43").Value = "  -1.56%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "22.08"
$ws.Range("E44").Value = "  +0.35%  "

# Rows 49 and 50 swap their Coin/Link content (and get new Price/Volume values)
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "5.43"
$ws.Range("E49").Value = "  -5.07%  "

$ws.Range("B50").Value = "SEI"
$ws.Range("C50").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D50").Value = "0.906"
$ws.Range("E50").Value = "  -2.56%  "

# Remaining price/volume updates for all other rows
$ws.Range("D2").Value = "51.642.76"
$ws.Range("E2").Value = "  -0.59%  "

$ws.Range("D3").Value = "2.782.95"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "352.17"
$ws.Range("E5").Value = "  -1.64%  "

$ws.Range("D6").Value = "108.85"
$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("D7").Value = "0.550"
$ws.Range("E7").Value = "  -2.53%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +2.14%  "

$ws.Range("D10").Value = "39.81"
$ws.Range("E10").Value = "  -0.48%  "

$ws.Range("E11").Value = "  +2.64%  "

$ws.Range("D12").Value = "20.17"
$ws.Range("E12").Value = "  +3.69%  "

$ws.Range("E13").Value = "  -2.02%  "

$ws.Range("D14").Value = "7.67"
$ws.Range("E14").Value = "  +1.27%  "

$ws.Range("D15").Value = "3.221.94"
$ws.Range("E15").Value = "  -0.19%  "

$ws.Range("D16").Value = "2.793.99"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("E17").Value = "  -2.28%  "

$ws.Range("D18").Value = "51.618.77"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").Value = "7.74"
$ws.Range("E19").Value = "  +4.59%  "

$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("E21").Value = "  +1.42%  "

$ws.Range("D22").Value = "0.0₃0965"
$ws.Range("E22").Value = "  -1.62%  "

$ws.Range("D23").Value = "69.88"
$ws.Range("E23").Value = "  -0.53%  "

$ws.Range("D24").Value = "267.03"
$ws.Range("E24").Value = "  -2.66%  "

$ws.Range("D25").Value = "2.73"
$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("D26").Value = "26.12"
$ws.Range("E26").Value = "  -2.12%  "

$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.10%  "

$ws.Range("E28").Value = "  +12.28%  "

$ws.Range("D29").Value = "10.21"
$ws.Range("E29").Value = "  +0.39%  "

$ws.Range("D30").Value = "36.99"
$ws.Range("E30").Value = "  +7.34%  "

$ws.Range("D31").Value = "2.22"
$ws.Range("E31").Value = "  -2.19%  "

$ws.Range("E32").Value = "  +8.41%  "

$ws.Range("D33").Value = "51.72"
$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("D34").Value = "0.0453"
$ws.Range("E34").Value = "  -2.56%  "

$ws.Range("D35").Value = "5.57"
$ws.Range("E35").Value = "  +5.07%  "

$ws.Range("D36").Value = "0.0831"
$ws.Range("E36").Value = "  -1.83%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").Value = "18.48"
$ws.Range("E38").Value = "  +2.58%  "

$ws.Range("E39").Value = "  -2.38%  "

$ws.Range("E40").Value = "  -1.52%  "

$ws.Range("E41").Value = "  -1.13%  "

$ws.Range("E42").Value = "  -0.61%  "

$ws.Range("E45").Value = "  -2.71%  "

$ws.Range("D46").Value = "2.121.07"
$ws.Range("E46").Value = "  +2.25%  "

$ws.Range("D47").Value = "3.30"
$ws.Range("E47").Value = "  +1.44%  "

$ws.Range("E48").Value = "  +6.54%  "

$ws.Range("E51").Value = "  +9.08%  "
